$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3319.2974
$ws.Range("J17").Value = 3393.0833
$ws.Range("L17").Value = 10179.2499
$ws.Range("N17").Value = -10515.2499
$ws.Range("H19").Value = 731
$ws.Range("I19").Value = 523.5
$ws.Range("K19").Value = 523.5
$ws.Range("M19").Value = -348.5
$ws.Range("H43").Value = 1626.2
$ws.Range("J43").Value = 1596.5
$ws.Range("L43").Value = 1596.5
$ws.Range("N43").Value = -1734.5
$ws.Range("H80").Value = 2854484
$ws.Range("I80").Value = 1357.2858
$ws.Range("J80").Value = 3853078.2
$ws.Range("K80").Value = 4071.8574
$ws.Range("L80").Value = 11559234.6
$ws.Range("M80").Value = -3073.8574
$ws.Range("N80").Value = -11561230.6
$ws.Range("H83").Value = 2854484
$ws.Range("I83").Value = 1357.2858
$ws.Range("J83").Value = 3853078.2
$ws.Range("K83").Value = 12215.5722
$ws.Range("L83").Value = 34677703.8
$ws.Range("M83").Value = -7223.572200000001
$ws.Range("N83").Value = -34687687.8
$ws.Range("H86").Value = 6648.3687
$ws.Range("I86").Value = 1440.2727
$ws.Range("K86").Value = 1440.2727
$ws.Range("M86").Value = -317.2727
$ws.Range("H89").Value = 6648.3687
$ws.Range("I89").Value = 1440.2727
$ws.Range("K89").Value = 7201.363499999999
$ws.Range("M89").Value = -1585.363499999999
$ws.Range("H116").Value = 4344.375
$ws.Range("I116").Value = 1750.625
$ws.Range("J116").Value = 6938.125
$ws.Range("K116").Value = 1750.625
$ws.Range("L116").Value = 6938.125
$ws.Range("M116").Value = 1691.375
$ws.Range("N116").Value = -13822.125
$ws.Range("H129").Value = 846.0833
$ws.Range("J129").Value = 850.2727
$ws.Range("L129").Value = 2550.8181
$ws.Range("N129").Value = -12550.8181
$ws.Range("H135").Value = 26322544
$ws.Range("I135").Value = 809.3333
$ws.Range("K135").Value = 7283.9997
$ws.Range("M135").Value = -4748.9997
$ws.Range("H138").Value = 1653.0769
$ws.Range("I138").Value = 528
$ws.Range("K138").Value = 1584
$ws.Range("M138").Value = 3556

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 104206.1
$ws.Range("I88").Value = 2000
$ws.Range("J88").Value = 129757.625
$ws.Range("K88").Value = 2000
$ws.Range("L88").Value = 129757.625
$ws.Range("M88").Value = -1594
$ws.Range("N88").Value = -130569.625
$ws.Range("H91").Value = 104206.1
$ws.Range("I91").Value = 2000
$ws.Range("J91").Value = 129757.625
$ws.Range("K91").Value = 2000
$ws.Range("L91").Value = 129757.625
$ws.Range("M91").Value = -596
$ws.Range("N91").Value = -132565.625
$ws.Range("H102").Value = 3562.7273
$ws.Range("I102").Value = 1330
$ws.Range("J102").Value = 7470
$ws.Range("K102").Value = 1330
$ws.Range("L102").Value = 7470
$ws.Range("M102").Value = 292
$ws.Range("N102").Value = -10714
$ws.Range("H132").Value = 25701.955
$ws.Range("I132").Value = 2580.5386
$ws.Range("K132").Value = 7741.6158
$ws.Range("M132").Value = -5211.6158

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2063.9
$ws.Range("I86").Value = 1854.875
$ws.Range("K86").Value = 1854.875
$ws.Range("M86").Value = -731.875
$ws.Range("H89").Value = 2063.9
$ws.Range("I89").Value = 1854.875
$ws.Range("K89").Value = 9274.375
$ws.Range("M89").Value = -3658.375
$ws.Range("H94").Value = 2592.2432
$ws.Range("I94").Value = 1372.52
$ws.Range("K94").Value = 1372.52
$ws.Range("M94").Value = -921.52
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 51000
$ws.Range("J138").Value = 51000
$ws.Range("L138").Value = 51000
$ws.Range("N138").Value = -61280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 15316.667
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 15316.667
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 15316.667
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -16566.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 1000
$ws.Range("J19").Value = 1000
$ws.Range("L19").Value = 3000
$ws.Range("N19").Value = -3348
$ws.Range("H113").Value = 21381.3
$ws.Range("I113").Value = 66810
$ws.Range("J113").Value = 1911.8572
$ws.Range("K113").Value = 200430
$ws.Range("L113").Value = 5735.571599999999
$ws.Range("M113").Value = -198260
$ws.Range("N113").Value = -10075.5716
$ws.Range("H131").Value = 771.35
$ws.Range("I131").Value = 364.75
$ws.Range("J131").Value = 788.2917
$ws.Range("K131").Value = 1094.25
$ws.Range("L131").Value = 2364.8751
$ws.Range("M131").Value = 3945.75
$ws.Range("N131").Value = -12444.8751
$ws.Range("H132").Value = 922.1539
$ws.Range("I132").Value = 947.2727
$ws.Range("J132").Value = 784
$ws.Range("K132").Value = 8525.454299999999
$ws.Range("L132").Value = 7056
$ws.Range("M132").Value = -5995.454299999999
$ws.Range("N132").Value = -12116

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2508179.8
$ws.Range("I70").Value = 12808.818
$ws.Range("K70").Value = 12808.818
$ws.Range("M70").Value = -12538.818
$ws.Range("H73").Value = 2508179.8
$ws.Range("I73").Value = 12808.818
$ws.Range("K73").Value = 12808.818
$ws.Range("M73").Value = -11872.818
$ws.Range("H126").Value = 4823.793
$ws.Range("I126").Value = 3336.875
$ws.Range("J126").Value = 6653.846
$ws.Range("K126").Value = 10010.625
$ws.Range("L126").Value = 19961.538
$ws.Range("M126").Value = -7540.625
$ws.Range("N126").Value = -24901.538
$ws.Range("H132").Value = 255333.17
$ws.Range("J132").Value = 128499.5
$ws.Range("L132").Value = 385498.5
$ws.Range("N132").Value = -390558.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2196.8333
$ws.Range("I22").Value = 5200.5
$ws.Range("J22").Value = 695
$ws.Range("K22").Value = 5200.5
$ws.Range("L22").Value = 695
$ws.Range("M22").Value = -4905.5
$ws.Range("N22").Value = -1285
$ws.Range("H27").Value = 2196.8333
$ws.Range("I27").Value = 5200.5
$ws.Range("J27").Value = 695
$ws.Range("K27").Value = 5200.5
$ws.Range("L27").Value = 695
$ws.Range("M27").Value = -5093.5
$ws.Range("N27").Value = -909
$ws.Range("H93").Value = 2700.6155
$ws.Range("I93").Value = 2828
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 2828
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = -1580
$ws.Range("N93").Value = -4496
$ws.Range("H132").Value = 2288.1875
$ws.Range("I132").Value = 1642.3334
$ws.Range("J132").Value = 4225.75
$ws.Range("K132").Value = 4927.0002
$ws.Range("L132").Value = 12677.25
$ws.Range("M132").Value = -2397.0002
$ws.Range("N132").Value = -17737.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 775
$ws.Range("I100").Value = 466.66666
$ws.Range("J100").Value = 1700
$ws.Range("K100").Value = 933.33332
$ws.Range("L100").Value = 3400
$ws.Range("M100").Value = -392.33332
$ws.Range("N100").Value = -4482
$ws.Range("H132").Value = 4499.5
$ws.Range("I132").Value = 4000
$ws.Range("K132").Value = 12000
$ws.Range("M132").Value = -9470
